# Update the Feb 13th progress-log row: note that the build now compiles
# on all toolchains (CLion/cmake/clang), and move the active selection
# to C14 to reflect where editing continued next.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C9").Value = "Token Stack complete; Compiles on all"

$ws.Range("C14").Select()
